$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Fgf7"
$ws.Cells.Item(2, 3).Value = "Fgfr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1658776666666667
$ws.Cells.Item(2, 8).Value = 0.497633
$ws.Cells.Item(2, 9).Value = 0.01379803068909966
$ws.Cells.Item(2, 10).Value = 0.01379803068909966
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8155003333333334
$ws.Cells.Item(2, 14).Value = 2.446501
$ws.Cells.Item(2, 15).Value = 0.1910612426590028
$ws.Cells.Item(2, 16).Value = 0.1910612426590029
$ws.Cells.Item(2, 17).Value = 0.1352732924592222
$ws.Cells.Item(2, 18).Value = 1.217459632133
$ws.Cells.Item(2, 19).Value = 0.002636268889706438
$ws.Cells.Item(2, 20).Value = 0.002636268889706439

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Fgf7"
$ws.Cells.Item(3, 3).Value = "Fgfr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1658776666666667
$ws.Cells.Item(3, 8).Value = 0.497633
$ws.Cells.Item(3, 9).Value = 0.01379803068909966
$ws.Cells.Item(3, 10).Value = 0.01379803068909966
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.333134333333334
$ws.Cells.Item(3, 14).Value = 9.999403000000001
$ws.Cells.Item(3, 15).Value = 0.7809105179307759
$ws.Cells.Item(3, 16).Value = 0.780910517930776
$ws.Cells.Item(3, 17).Value = 0.552892545899889
$ws.Cells.Item(3, 18).Value = 4.976032913099
$ws.Cells.Item(3, 19).Value = 0.01077502729184955
$ws.Cells.Item(3, 20).Value = 0.01077502729184956

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Fgf7"
$ws.Cells.Item(4, 3).Value = "Fgfr2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1658776666666667
$ws.Cells.Item(4, 8).Value = 0.497633
$ws.Cells.Item(4, 9).Value = 0.01379803068909966
$ws.Cells.Item(4, 10).Value = 0.01379803068909966
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.119632
$ws.Cells.Item(4, 14).Value = 0.358896
$ws.Cells.Item(4, 15).Value = 0.02802823941022116
$ws.Cells.Item(4, 16).Value = 0.02802823941022117
$ws.Cells.Item(4, 17).Value = 0.01984427701866667
$ws.Cells.Item(4, 18).Value = 0.178598493168
$ws.Cells.Item(4, 19).Value = 0.0003867345075436641
$ws.Cells.Item(4, 20).Value = 0.0003867345075436642

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Fgf7"
$ws.Cells.Item(5, 3).Value = "Fgfr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 11.51839566666667
$ws.Cells.Item(5, 8).Value = 34.555187
$ws.Cells.Item(5, 9).Value = 0.9581228147923823
$ws.Cells.Item(5, 10).Value = 0.9581228147923824
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8155003333333334
$ws.Cells.Item(5, 14).Value = 2.446501
$ws.Cells.Item(5, 15).Value = 0.1910612426590028
$ws.Cells.Item(5, 16).Value = 0.1910612426590029
$ws.Cells.Item(5, 17).Value = 9.39325550563189
$ws.Cells.Item(5, 18).Value = 84.539299550687
$ws.Cells.Item(5, 19).Value = 0.1830601356141742
$ws.Cells.Item(5, 20).Value = 0.1830601356141742

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Fgf7"
$ws.Cells.Item(6, 3).Value = "Fgfr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 11.51839566666667
$ws.Cells.Item(6, 8).Value = 34.555187
$ws.Cells.Item(6, 9).Value = 0.9581228147923823
$ws.Cells.Item(6, 10).Value = 0.9581228147923824
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 3.333134333333334
$ws.Cells.Item(6, 14).Value = 9.999403000000001
$ws.Cells.Item(6, 15).Value = 0.7809105179307759
$ws.Cells.Item(6, 16).Value = 0.780910517930776
$ws.Cells.Item(6, 17).Value = 38.39236006148457
$ws.Cells.Item(6, 18).Value = 345.531240553361
$ws.Cells.Item(6, 19).Value = 0.7482081835408121
$ws.Cells.Item(6, 20).Value = 0.7482081835408123

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Fgf7"
$ws.Cells.Item(7, 3).Value = "Fgfr2"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 11.51839566666667
$ws.Cells.Item(7, 8).Value = 34.555187
$ws.Cells.Item(7, 9).Value = 0.9581228147923823
$ws.Cells.Item(7, 10).Value = 0.9581228147923824
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.119632
$ws.Cells.Item(7, 14).Value = 0.358896
$ws.Cells.Item(7, 15).Value = 0.02802823941022116
$ws.Cells.Item(7, 16).Value = 0.02802823941022117
$ws.Cells.Item(7, 17).Value = 1.377968710394667
$ws.Cells.Item(7, 18).Value = 12.401718393552
$ws.Cells.Item(7, 19).Value = 0.02685449563739588
$ws.Cells.Item(7, 20).Value = 0.02685449563739589

# Row 8
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Fgf7"
$ws.Cells.Item(8, 3).Value = "Fgfr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.1594243333333333
$ws.Cells.Item(8, 8).Value = 0.478273
$ws.Cells.Item(8, 9).Value = 0.01326122972505393
$ws.Cells.Item(8, 10).Value = 0.01326122972505393
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.8155003333333334
$ws.Cells.Item(8, 14).Value = 2.446501
$ws.Cells.Item(8, 15).Value = 0.1910612426590028
$ws.Cells.Item(8, 16).Value = 0.1910612426590029
$ws.Cells.Item(8, 17).Value = 0.1300105969747778
$ws.Cells.Item(8, 18).Value = 1.170095372773
$ws.Cells.Item(8, 19).Value = 0.00253370703045531
$ws.Cells.Item(8, 20).Value = 0.002533707030455311

# Row 9
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Fgf7"
$ws.Cells.Item(9, 3).Value = "Fgfr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.1594243333333333
$ws.Cells.Item(9, 8).Value = 0.478273
$ws.Cells.Item(9, 9).Value = 0.01326122972505393
$ws.Cells.Item(9, 10).Value = 0.01326122972505393
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 3.333134333333334
$ws.Cells.Item(9, 14).Value = 9.999403000000001
$ws.Cells.Item(9, 15).Value = 0.7809105179307759
$ws.Cells.Item(9, 16).Value = 0.780910517930776
$ws.Cells.Item(9, 17).Value = 0.5313827190021112
$ws.Cells.Item(9, 18).Value = 4.782444471019001
$ws.Cells.Item(9, 19).Value = 0.01035583377299086
$ws.Cells.Item(9, 20).Value = 0.01035583377299087

# Row 10
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Fgf7"
$ws.Cells.Item(10, 3).Value = "Fgfr2"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.1594243333333333
$ws.Cells.Item(10, 8).Value = 0.478273
$ws.Cells.Item(10, 9).Value = 0.01326122972505393
$ws.Cells.Item(10, 10).Value = 0.01326122972505393
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.119632
$ws.Cells.Item(10, 14).Value = 0.358896
$ws.Cells.Item(10, 15).Value = 0.02802823941022116
$ws.Cells.Item(10, 16).Value = 0.02802823941022117
$ws.Cells.Item(10, 17).Value = 0.01907225184533333
$ws.Cells.Item(10, 18).Value = 0.171650266608
$ws.Cells.Item(10, 19).Value = 0.0003716889216077528
$ws.Cells.Item(10, 20).Value = 0.0003716889216077529

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Fgf7"
$ws.Cells.Item(11, 3).Value = "Fgfr2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.1781386666666667
$ws.Cells.Item(11, 8).Value = 0.534416
$ws.Cells.Item(11, 9).Value = 0.01481792479346402
$ws.Cells.Item(11, 10).Value = 0.01481792479346403
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.8155003333333334
$ws.Cells.Item(11, 14).Value = 2.446501
$ws.Cells.Item(11, 15).Value = 0.1910612426590028
$ws.Cells.Item(11, 16).Value = 0.1910612426590029
$ws.Cells.Item(11, 17).Value = 0.1452721420462222
$ws.Cells.Item(11, 18).Value = 1.307449278416
$ws.Cells.Item(11, 19).Value = 0.002831131124666884
$ws.Cells.Item(11, 20).Value = 0.002831131124666885

# Row 12
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Fgf7"
$ws.Cells.Item(12, 3).Value = "Fgfr2"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.1781386666666667
$ws.Cells.Item(12, 8).Value = 0.534416
$ws.Cells.Item(12, 9).Value = 0.01481792479346402
$ws.Cells.Item(12, 10).Value = 0.01481792479346403
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 3.333134333333334
$ws.Cells.Item(12, 14).Value = 9.999403000000001
$ws.Cells.Item(12, 15).Value = 0.7809105179307759
$ws.Cells.Item(12, 16).Value = 0.780910517930776
$ws.Cells.Item(12, 17).Value = 0.5937601059608889
$ws.Cells.Item(12, 18).Value = 5.343840953648001
$ws.Cells.Item(12, 19).Value = 0.01157147332512328
$ws.Cells.Item(12, 20).Value = 0.01157147332512328

# Row 13
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Fgf7"
$ws.Cells.Item(13, 3).Value = "Fgfr2"
$ws.Cells.Item(13, 4).Value = "MuSCs"
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.1781386666666667
$ws.Cells.Item(13, 8).Value = 0.534416
$ws.Cells.Item(13, 9).Value = 0.01481792479346402
$ws.Cells.Item(13, 10).Value = 0.01481792479346403
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.119632
$ws.Cells.Item(13, 14).Value = 0.358896
$ws.Cells.Item(13, 15).Value = 0.02802823941022116
$ws.Cells.Item(13, 16).Value = 0.02802823941022117
$ws.Cells.Item(13, 17).Value = 0.02131108497066667
$ws.Cells.Item(13, 18).Value = 0.191799764736
$ws.Cells.Item(13, 19).Value = 0.0004153203436738616
$ws.Cells.Item(13, 20).Value = 0.0004153203436738617
